$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "error:X%" labels in column A (rows 3-8) to "diff:X%"
$ws.Range("A3").Value = "diff:0.00%"
$ws.Range("A4").Value = "diff:5.00%"
$ws.Range("A5").Value = "diff:10.00%"
$ws.Range("A6").Value = "diff:15.00%"
$ws.Range("A7").Value = "diff:20.00%"
$ws.Range("A8").Value = "diff:25.00%"
